$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '69.177.68'
$ws.Range("E2").Value = '  +1.14%  '

# Row 3
$ws.Range("D3").Value = '3.773.91'
$ws.Range("E3").Value = '  -0.68%  '

# Row 4
$ws.Range("E4").Value = '  -0.39%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '633.17'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +4.06%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '166.46'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +1.76%  '

# Row 7
$ws.Range("D7").Value = '3.771.58'
$ws.Range("E7").Value = '  -0.63%  '

# Row 8
$ws.Range("E8").Value = '  -0.10%  '

# Row 9
$ws.Range("E9").Value = '  +0.92%  '

# Row 10
$ws.Range("E10").Value = '  -0.25%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.461'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +2.42%  '

# Row 12
$ws.Range("E12").Value = '  -2.90%  '

# Row 13
$ws.Range("E13").Value = '  -1.40%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '34.95'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.41%  '

# Row 15
$ws.Range("D15").Value = '4.406.46'
$ws.Range("E15").Value = '  -0.69%  '

# Row 16
$ws.Range("D16").Value = '3.746.31'
$ws.Range("E16").Value = '  -1.40%  '

# Row 17
$ws.Range("D17").Value = '69.172.57'
$ws.Range("E17").Value = '  +1.06%  '

# Row 18
$ws.Range("E18").Value = '  -2.29%  '

# Row 19
$ws.Range("E19").Value = '  +0.68%  '

# Row 20
$ws.Range("E20").Value = '  -0.70%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '463.78'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.39%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '9.56'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.25%  '

# Row 23
$ws.Range("E23").Value = '  +1.50%  '

# Row 24
$ws.Range("E24").Value = '  -0.66%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '82.78'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.77%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '12.11'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.90%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.15'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +2.03%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.12'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +1.39%  '

# Row 29
$ws.Range("E29").Value = '  +0.00%  '

# Row 30
$ws.Range("D30").Value = '3.923.11'
$ws.Range("E30").Value = '  -0.62%  '

# Row 31
$ws.Range("E31").Value = '  +2.09%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.28'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +3.33%  '

# Row 33
$ws.Range("E33").Value = '  -1.66%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '28.53'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -1.72%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.168'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +15.02%  '

# Row 36
$ws.Range("E36").Value = '  -0.17%  '

# Row 37
$ws.Range("B37").Value = 'Aptos'
$ws.Range("C37").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '8.99'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.65%  '

# Row 38
$ws.Range("B38").Value = 'RenzoRestakedETH'
$ws.Range("C38").Value = 'https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth'
$ws.Range("D38").Value = '3.726.03'
$ws.Range("E38").Value = '  -0.55%  '

# Row 39
$ws.Range("E39").Value = '  +0.46%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.34'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +4.61%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '5.81'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -1.07%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.963'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -1.54%  '

# Row 43
$ws.Range("E43").Value = '  -0.08%  '

# Row 44
$ws.Range("E44").Value = '  -0.01%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '158.32'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +3.39%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.97'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +5.84%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.43'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +1.19%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '43.25'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.54%  '

# Row 49
$ws.Range("E49").Value = '  -0.31%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '46.75'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.25%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '8.39'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.24%  '
